$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number-looking string (e.g. '1.002').
# These must be forced to text so Excel doesn't silently convert them to
# numeric values (the source data stores prices as literal text).
$textForceCells = [ordered]@{
    'D4' = '1.002'
    'D5' = '325.80'
    'D7' = '0.4632'
    'D9' = '0.07871'
    'D10' = '0.9590'
    'D11' = '21.83'
    'D13' = '5.668'
    'D14' = '6.898'
    'D15' = '0.06776'
    'D16' = '87.15'
    'D17' = '1.002'
    'D18' = '0.000009926'
    'D19' = '16.61'
    'D20' = '1.001'
    'D22' = '5.313'
    'D23' = '10.97'
    'D24' = '2.090'
    'D26' = '153.78'
    'D27' = '19.16'
    'D28' = '5.734'
    'D29' = '1.975'
    'D30' = '117.26'
    'D31' = '0.9381'
    'D32' = '0.09244'
    'D33' = '5.294'
    'D34' = '1.317'
    'D35' = '3.285'
    'D36' = '0.05863'
    'D37' = '0.02143'
    'D38' = '1.148'
    'D39' = '7.760'
    'D40' = '0.5579'
    'D41' = '9.873'
    'D42' = '0.1761'
    'D43' = '11.60'
    'D44' = '0.5272'
    'D45' = '0.07001'
    'D46' = '1.131'
    'D47' = '2.121'
    'D48' = '1.833'
    'D49' = '112.95'
    'D50' = '1.001'
    'D51' = '2.320'
}

foreach ($addr in $textForceCells.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textForceCells[$addr]
}

# Remaining cells (text that Excel will not misinterpret as numbers,
# e.g. percentages with surrounding spaces, multi-dot prices, coin names,
# and links) can be written directly.
$directCells = [ordered]@{
    'D2' = '28.015.68'
    'E2' = '  -1.98%  '
    'D3' = '1.830.29'
    'E3' = '  -1.02%  '
    'E4' = '  -0.10%  '
    'E5' = '  -2.94%  '
    'E6' = '  -0.08%  '
    'E8' = '  -0.97%  '
    'E9' = '  -0.14%  '
    'E10' = '  -2.31%  '
    'E11' = '  -1.53%  '
    'D12' = '1.797.19'
    'E12' = '  -4.37%  '
    'E13' = '  -2.98%  '
    'E14' = '  -1.50%  '
    'E15' = '  -0.85%  '
    'E16' = '  -0.52%  '
    'E17' = '  -0.13%  '
    'E18' = '  -1.81%  '
    'E19' = '  -2.31%  '
    'E20' = '  -0.07%  '
    'D21' = '28.041.63'
    'E21' = '  -1.96%  '
    'E22' = '  -1.46%  '
    'E23' = '  -2.40%  '
    'E24' = '  -1.58%  '
    'D25' = '2.043.27'
    'E25' = '  -4.41%  '
    'E26' = '  +0.38%  '
    'E27' = '  -1.15%  '
    'E28' = '  -7.37%  '
    'E29' = '  -2.15%  '
    'E30' = '  -0.11%  '
    'E31' = '  -3.73%  '
    'E32' = '  -2.15%  '
    'E33' = '  -1.37%  '
    'E34' = '  -2.33%  '
    'E35' = '  -6.22%  '
    'E36' = '  -4.59%  '
    'E37' = '  -2.12%  '
    'E38' = '  -1.13%  '
    'E40' = '  -1.89%  '
    'E41' = '  -2.25%  '
    'E42' = '  -1.55%  '
    'E43' = '  -2.45%  '
    'E44' = '  -1.94%  '
    'E45' = '  -1.88%  '
    'B46' = 'WEMIXToken'
    'C46' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'E46' = '  -9.48%  '
    'B47' = 'RenderToken'
    'C47' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'E47' = '  -10.82%  '
    'B48' = 'NEARProtocol'
    'C48' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E48' = '  -3.70%  '
    'E49' = '  -0.14%  '
    'E50' = '  -0.09%  '
    'E51' = '  +0.54%  '
}

foreach ($addr in $directCells.Keys) {
    $ws.Range($addr).Value = $directCells[$addr]
}
